$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find a top-level paragraph whose visible text (paragraph mark
# stripped) equals $text exactly, and return it.
# ---------------------------------------------------------------------------
function Get-ParagraphByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. Strike through the five "cut" timeline bullets (paragraph mark + run
#    both get <w:strike/>, matching how Word records a whole-paragraph
#    strikethrough toggle).
# ---------------------------------------------------------------------------
$strikeTexts = @(
    "Jack is unhappy in life and sees ad on television to join the RRH",
    "Upon joining, meets Nicholas, another like him",
    "Receive grueling training on survival in the wilderness",
    "Jack runs towards forest",
    "Jack runs into great creature, kills with a knife"
)

foreach ($txt in $strikeTexts) {
    $p = Get-ParagraphByText $d $txt
    if ($p -ne $null) {
        $p.Range.Font.StrikeThrough = 1
    }
}

# ---------------------------------------------------------------------------
# 2. Append a new sentence, as its own run, right after "Nick confronts him".
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Nick confronts him")
$r.Collapse(0)
$insertStart = $r.Start
$addition = " (at some point, there is a knife to jack’s throat.)"
$r.InsertAfter($addition)

# Force the newly inserted text to live in its own <w:r> (rather than being
# silently coalesced back into the preceding "Nick confronts him" run) by
# toggling a character property on just that span and restoring it.
$newRange = $d.Range($insertStart, $insertStart + $addition.Length)
$newRange.Font.Bold = 1
$newRange.Font.Bold = 0
